$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 4 and row 5 for columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $addr4 = "${col}4"
    $addr5 = "${col}5"
    $cell4 = $ws.Range($addr4)
    $cell5 = $ws.Range($addr5)
    $val4 = $cell4.Value2
    $val5 = $cell5.Value2
    $cell4.Value2 = $val5
    $cell5.Value2 = $val4
}
